$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.472.96'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = '3.677.52'
$ws.Range("E3").Value = '  -3.45%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.28'
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.49'
$ws.Range("E6").Value = '  -5.17%  '
$ws.Range("D7").Value = '3.678.71'
$ws.Range("E7").Value = '  -3.23%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("E9").Value = '  +0.87%  '
$ws.Range("E10").Value = '  +2.72%  '
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("E12").Value = '  -1.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.80'
$ws.Range("E13").Value = '  -1.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000243'
$ws.Range("E14").Value = '  -0.94%  '
$ws.Range("D15").Value = '4.291.66'
$ws.Range("E15").Value = '  -3.26%  '
$ws.Range("D16").Value = '3.673.85'
$ws.Range("E16").Value = '  -3.32%  '
$ws.Range("D17").Value = '68.487.63'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.20'
$ws.Range("E18").Value = '  +0.63%  '
$ws.Range("E19").Value = '  -0.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.04'
$ws.Range("E20").Value = '  +4.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '491.76'
$ws.Range("E21").Value = '  +0.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.07'
$ws.Range("E22").Value = '  -2.32%  '
$ws.Range("E23").Value = '  -2.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.36'
$ws.Range("E24").Value = '  -0.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000140'
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("E26").Value = '  -5.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.14'
$ws.Range("E27").Value = '  -1.17%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.08'
$ws.Range("E28").Value = '  -1.57%  '
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.91'
$ws.Range("E30").Value = '  -0.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.79'
$ws.Range("E31").Value = '  +0.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.36'
$ws.Range("E32").Value = '  -3.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.22'
$ws.Range("E33").Value = '  -4.77%  '
$ws.Range("D34").Value = '3.813.37'
$ws.Range("E34").Value = '  -3.28%  '
$ws.Range("E35").Value = '  -1.28%  '
$ws.Range("D36").Value = '3.617.94'
$ws.Range("E36").Value = '  -3.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.70'
$ws.Range("E39").Value = '  -1.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.132'
$ws.Range("E40").Value = '  -4.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.320'
$ws.Range("E41").Value = '  -2.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '48.93'
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '427.05'
$ws.Range("E43").Value = '  -5.61%  '
$ws.Range("E44").Value = '  -2.84%  '
$ws.Range("E45").Value = '  -3.73%  '
$ws.Range("E46").Value = '  +0.67%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.03'
$ws.Range("E48").Value = '  -3.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.23'
$ws.Range("E49").Value = '  +1.54%  '
$ws.Range("E50").Value = '  -1.47%  '
$ws.Range("D51").Value = '2.717.62'
$ws.Range("E51").Value = '  -3.89%  '
